# "did multiple rtk query working"
# Insert a new "debt" worksheet between "Download" and "Delete", populate it
# with a short backlog/task list in column C, size the column to fit, and
# leave the new sheet selected/active (G17 selected) as the final state.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "Download" (i.e. before "Delete"),
# so the final tab order is Download, debt, Delete, Tech_Data_Flow.
$downloadSheet = $wb.Worksheets.Item("Download")
$debt = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $downloadSheet)
$debt.Name = "debt"

# Task / backlog entries in column C.
$debt.Range("C2").Value = "delete with filter"
$debt.Range("C3").Value = "encrytion of payload"
$debt.Range("C4").Value = "refactor all api to HTTPS"
$debt.Range("C5").Value = "Fill patern in a cell"
$debt.Range("C8").Value = "https://sproutsocial.com/insights/linkedin-business-page/"

# Widen column C so the text fits.
$debt.Columns.Item(3).ColumnWidth = 71

# Leave this as the selected cell/sheet, matching the saved view state.
$debt.Range("G17").Select()
